$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questions Tracker")

$ws.Range("B35").Value = "Pivot in an Array"
$ws.Range("D35").Value = 1
$ws.Range("E35").Value = "15/11/2022"
$ws.Range("F35").Value = 2
$ws.Range("G35").Value = "Yes"
$ws.Range("H35").Value = "Good Question"

$ws.Range("I35").Select()
